# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.401.22"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "3.088.03"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.45"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.71"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.086.92"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.17"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "3.617.04"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.83"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "57.468.68"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "3.091.66"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.75"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.02"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "338.26"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.510"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.47"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "0.0₃0907"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("B29").Value = "USDe"
$ws.Range("C29").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.46"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.15"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.86"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.32"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.99"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0657"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.53"
$ws.Range("E41").Value = "  +11.09%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.93"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("B43").Value = "RenzoRestakedETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D43").Value = "3.129.38"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.683"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.66"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "2.301.21"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0258"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.59"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.01"
$ws.Range("E51").Value = "  +1.43%  "
